$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above current row 2; the existing "NN" row (with all
# its values) shifts down to row 4.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# The freshly inserted rows 2 and 3 pick up a bold/no-border style copied
# down from the header row; strip that so the new rows start unstyled like
# the rest of the numeric data cells.
$ws.Rows.Item(2).ClearFormats()
$ws.Rows.Item(3).ClearFormats()

# Restore the bordered/bold/centered label style on A2 and A3 by copying the
# format that already sits on A4 (the shifted-down original "NN" label).
$ws.Range("A4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: "mf" ---
$ws.Range("A2").Value = "mf"
$ws.Range("B2").Value = 1.032138758722124
$ws.Range("C2").Value = 0.7652540722461235
$ws.Range("D2").Value = 0.01731744981248621
$ws.Range("E2").Value = 0.180425360411634
$ws.Range("F2").Value = 0.04113263785394933
$ws.Range("G2").Value = 0.003100835403920611
$ws.Range("H2").Value = 0.005311706154370903
$ws.Range("I2").Value = 0.04087925014736316
$ws.Range("J2").Value = 0.001466035986383272
$ws.Range("K2").Value = 0.08211624441132637
$ws.Range("L2").Value = 0.00192577578732003
$ws.Range("M2").Value = 0.01731744981248621
# N2 stays an explicit empty text cell (matches column's "no value yet" marker).
$ws.Range("N2").Value = "'"
$ws.Range("N2").ClearFormats()

# --- Row 3: "mmr" ---
$ws.Range("A3").Value = "mmr"
$ws.Range("B3").Value = 1.009170831239415
$ws.Range("C3").Value = 0.7302371529881306
$ws.Range("D3").Value = 0.02172953893668652
$ws.Range("E3").Value = 0.1931026812795729
$ws.Range("F3").Value = 0.02175856929955291
$ws.Range("G3").Value = 0.00129118620221949
$ws.Range("H3").Value = 0.002334409549589788
$ws.Range("I3").Value = 0.02036073997695226
$ws.Range("J3").Value = 0.0005748124404458528
$ws.Range("K3").Value = 0.03899652260307997
$ws.Range("L3").Value = 0.0007728069700867142
$ws.Range("M3").Value = 0.02172953893668652
$ws.Range("N3").Value = "'"
$ws.Range("N3").ClearFormats()

# --- Row 4 (was row 2, "NN"): RMSE/MAE are no longer populated for this row ---
$ws.Range("B4").Value = "'"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = "'"
$ws.Range("C4").ClearFormats()
